$p = $ppt.ActivePresentation

# --- Update the cached "last saved" date/time fields -------------------
# These live on the Notes Master, the Slide Master, and every Slide
# Layout's "Date Placeholder" shape (shape index 3 in every layout / the
# master). Re-opening & re-saving the deck on 9/27/2023 refreshed the
# cached text of each field from 10/26/2022 to 9/27/2023.
$newDate = "9/27/2023"

$nm = $p.NotesMaster
$nmDateShape = $nm.Shapes.Item(2)
$nmDateShape.TextFrame.TextRange.Text = $newDate

$master = $p.Slides.Item(1).Master
$masterDateShape = $master.Shapes.Item(3)
$masterDateShape.TextFrame.TextRange.Text = $newDate

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    $layoutDateShape = $layout.Shapes.Item(3)
    $layoutDateShape.TextFrame.TextRange.Text = $newDate
}

# --- Remove the "Fall 2022" line from the title slide's subtitle -------
$slide1 = $p.Slides.Item(1)
$subtitle = $slide1.Shapes.Item(4)
$subtitle.TextFrame.TextRange.Text = "University of Mount Union"
